$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clean up the header row: the original header strings carried stray
# leading/trailing spaces (" ErrorDate", "EmployeeErrorId ", etc.) -
# rewrite them trimmed, in natural left-to-right column order.
$ws.Range("A1").Value = "EmployeeErrorId"
$ws.Range("B1").Value = "EmployeeId"
$ws.Range("C1").Value = "ErrorCodeId"
$ws.Range("D1").Value = "ErrorDate"

# Move/save the current selection to D7 (was C22).
$ws.Range("D7").Select()
